$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.028.18"
$ws.Range("E2").Value = "  -1.42%  "
$ws.Range("D3").Value = "2.059.76"
$ws.Range("E3").Value = "  +6.56%  "
$ws.Range("E4").Value = "  +0.18%  "
$c = $ws.Range("D5")
$c.Value = "'249.42"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -0.49%  "
$c = $ws.Range("D6")
$c.Value = "'0.650"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -5.74%  "
$c = $ws.Range("D7")
$c.Value = "'0.998"
$c.Style = "Normal"
$ws.Range("E7").Value = "  -0.19%  "
$c = $ws.Range("D8")
$c.Value = "'49.94"
$c.Style = "Normal"
$ws.Range("E8").Value = "  +3.99%  "
$c = $ws.Range("D9")
$c.Value = "'59.95"
$c.Style = "Normal"
$ws.Range("E9").Value = "  +2.01%  "
$c = $ws.Range("D10")
$c.Value = "'0.367"
$c.Style = "Normal"
$ws.Range("E10").Value = "  -3.04%  "
$c = $ws.Range("D11")
$c.Value = "'0.0734"
$c.Style = "Normal"
$ws.Range("E11").Value = "  -4.26%  "
$c = $ws.Range("D12")
$c.Value = "'0.105"
$c.Style = "Normal"
$ws.Range("E12").Value = "  +4.16%  "
$c = $ws.Range("D13")
$c.Value = "'14.91"
$c.Style = "Normal"
$ws.Range("E13").Value = "  -3.94%  "
$ws.Range("D14").Value = "2.365.64"
$ws.Range("E14").Value = "  +6.86%  "
$c = $ws.Range("D15")
$c.Value = "'0.823"
$c.Style = "Normal"
$ws.Range("E15").Value = "  -1.27%  "
$ws.Range("D16").Value = "2.107.08"
$ws.Range("E16").Value = "  +9.09%  "
$c = $ws.Range("D17")
$c.Value = "'5.04"
$c.Style = "Normal"
$ws.Range("E17").Value = "  -2.30%  "
$ws.Range("D18").Value = "36.891.28"
$ws.Range("E18").Value = "  -1.60%  "
$c = $ws.Range("D19")
$c.Value = "'71.69"
$c.Style = "Normal"
$ws.Range("E19").Value = "  -4.64%  "
$ws.Range("D20").Value = "0.0₃0816"
$ws.Range("E20").Value = "  -5.35%  "
$c = $ws.Range("D21")
$c.Value = "'13.11"
$c.Style = "Normal"
$ws.Range("E21").Value = "  -5.00%  "
$c = $ws.Range("D22")
$c.Value = "'237.67"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -5.99%  "
$c = $ws.Range("D23")
$c.Value = "'5.13"
$c.Style = "Normal"
$ws.Range("E23").Value = "  -1.58%  "
$ws.Range("E24").Value = "  -0.06%  "
$ws.Range("E25").Value = "  -2.12%  "
$c = $ws.Range("D26")
$c.Value = "'168.06"
$c.Style = "Normal"
$ws.Range("E26").Value = "  -0.56%  "
$c = $ws.Range("D27")
$c.Value = "'9.20"
$c.Style = "Normal"
$ws.Range("E27").Value = "  +3.56%  "
$c = $ws.Range("D28")
$c.Value = "'20.66"
$c.Style = "Normal"
$ws.Range("E28").Value = "  +9.65%  "
$c = $ws.Range("D29")
$c.Value = "'2.00"
$c.Style = "Normal"
$ws.Range("E29").Value = "  -6.30%  "
$c = $ws.Range("D30")
$c.Value = "'0.121"
$c.Style = "Normal"
$ws.Range("E30").Value = "  -5.82%  "
$c = $ws.Range("D31")
$c.Value = "'23.58"
$c.Style = "Normal"
$ws.Range("E31").Value = "  +25.31%  "
$c = $ws.Range("D32")
$c.Value = "'1.07"
$c.Style = "Normal"
$ws.Range("E32").Value = "  +18.27%  "
$c = $ws.Range("D33")
$c.Value = "'4.44"
$c.Style = "Normal"
$ws.Range("E33").Value = "  -3.37%  "
$c = $ws.Range("D34")
$c.Value = "'0.0600"
$c.Style = "Normal"
$ws.Range("E34").Value = "  -2.66%  "
$c = $ws.Range("D35")
$c.Value = "'0.0906"
$c.Style = "Normal"
$ws.Range("E35").Value = "  -1.70%  "
$c = $ws.Range("D36")
$c.Value = "'0.998"
$c.Style = "Normal"
$ws.Range("E36").Value = "  -0.28%  "
$ws.Range("E37").Value = "  +15.31%  "
$ws.Range("E38").Value = "  -2.85%  "
$c = $ws.Range("D39")
$c.Value = "'4.03"
$c.Style = "Normal"
$ws.Range("E39").Value = "  -7.00%  "
$ws.Range("E40").Value = "  -9.53%  "
$c = $ws.Range("D41")
$c.Value = "'17.49"
$c.Style = "Normal"
$ws.Range("E41").Value = "  -0.73%  "
$c = $ws.Range("D42")
$c.Value = "'0.0221"
$c.Style = "Normal"
$ws.Range("E42").Value = "  -3.03%  "
$ws.Range("E43").Value = "  +2.66%  "
$c = $ws.Range("D44")
$c.Value = "'96.71"
$c.Style = "Normal"
$ws.Range("E44").Value = "  -8.72%  "
$c = $ws.Range("D45")
$c.Value = "'2.76"
$c.Style = "Normal"
$ws.Range("E45").Value = "  -5.63%  "
$c = $ws.Range("D46")
$c.Value = "'0.0866"
$c.Style = "Normal"
$ws.Range("E46").Value = "  +2.95%  "
$c = $ws.Range("D47")
$c.Value = "'2.95"
$c.Style = "Normal"
$ws.Range("E47").Value = "  +4.84%  "
$ws.Range("D48").Value = "1.296.02"
$ws.Range("E48").Value = "  -3.83%  "
$ws.Range("B49").Value = "RocketPoolETH"
$ws.Range("C49").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D49").Value = "2.256.87"
$ws.Range("E49").Value = "  +7.00%  "
$ws.Range("B50").Value = "FraxShare"
$ws.Range("C50").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$c = $ws.Range("D50")
$c.Value = "'6.78"
$c.Style = "Normal"
$ws.Range("E50").Value = "  +5.47%  "
$c = $ws.Range("D51")
$c.Value = "'2.23"
$c.Style = "Normal"
$ws.Range("E51").Value = "  -7.40%  "
